# Apply English-language header translation to the "returned product" import
# template, and switch the active tab from "Produk Retur" to "Problem Produk".

$wb = $excel.ActiveWorkbook

$wsProdukRetur  = $wb.Worksheets.Item("Produk Retur")
$wsProblemProduk = $wb.Worksheets.Item("Problem Produk")

# ---------------------------------------------------------------------------
# Sheet "Produk Retur" (sheet1) -- translate table header row (row 1)
# ---------------------------------------------------------------------------
$wsProdukRetur.Range("A1").Value = "Serial Number"
$wsProdukRetur.Range("B1").Value = "Product Name"
$wsProdukRetur.Range("C1").Value = "Product Type"
$wsProdukRetur.Range("D1").Value = "Qty"
$wsProdukRetur.Range("E1").Value = "Customer (Optional)"

# ---------------------------------------------------------------------------
# Sheet "Problem Produk" (sheet2) -- translate table header row (row 1)
# ---------------------------------------------------------------------------
$wsProblemProduk.Range("A1").Value = "Serial Number"
$wsProblemProduk.Range("B1").Value = "Problem Component Name"
$wsProblemProduk.Range("C1").Value = "Note"
$wsProblemProduk.Range("D1").Value = "Status"

# Column A on "Problem Produk" is best-fit to its header text; now that the
# header is "Serial Number" instead of "Seri Produk" the best-fit width grows.
$wsProblemProduk.Columns.Item(1).ColumnWidth = 13

# ---------------------------------------------------------------------------
# Window / selection state: "Produk Retur" loses the active tab, and
# "Problem Produk" becomes the selected / active sheet.
# ---------------------------------------------------------------------------
[void]$wsProdukRetur.Activate()
[void]$wsProdukRetur.Range("B7").Select()

[void]$wsProblemProduk.Activate()
[void]$wsProblemProduk.Range("A2").Select()
